# This script inserts two new quarterly-reporting columns (D and E) into the
# PTGX financials worksheet, shifting the existing eight quarters of data to
# the right (old D:K -> new F:M), and populates the two new columns with the
# latest reported quarters (period ending 2018-12-31 and 2018-09-30), along
# with a handful of corrected values in the carried-over columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column D; this shifts the existing
# quarterly columns (D:K) two positions to the right (F:M).
$ws.Columns("D:E").Insert()

# The newly inserted columns have no formatting yet. Copy the number
# formats/styles from column F (the first of the shifted, still-formatted
# columns) down into D:E so the new cells match the existing look
# (date format in row 7/38/80, number format elsewhere).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns (and fix a few carried-over cells) with the
# updated financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2400
$ws.Range("E8").Value = 6100
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = 14200
$ws.Range("E12").Value = 12100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 17800
$ws.Range("E17").Value = 15500
$ws.Range("D18").Value = -15400
$ws.Range("E18").Value = -9400
$ws.Range("D20").Value = 700
$ws.Range("E20").Value = 700
$ws.Range("D21").Value = -14500
$ws.Range("E21").Value = -8600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -14700
$ws.Range("E23").Value = -8700
$ws.Range("D24").Value = -800
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "NA"
$ws.Range("I24").Value = "NA"
$ws.Range("J24").Value = "NA"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -13900
$ws.Range("E26").Value = -8700
$ws.Range("D27").Value = -13900
$ws.Range("E27").Value = -8700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -700
$ws.Range("E32").Value = -700
$ws.Range("D33").Value = -13900
$ws.Range("E33").Value = -8700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -13900
$ws.Range("E35").Value = -8700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 82200
$ws.Range("E41").Value = 87400
$ws.Range("D42").Value = 46600
$ws.Range("E42").Value = 51100
$ws.Range("D43").Value = 6000
$ws.Range("E43").Value = 5700
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 2600
$ws.Range("E45").Value = 3200
$ws.Range("D46").Value = 137500
$ws.Range("E46").Value = 147500
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 900
$ws.Range("E48").Value = 900
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1100
$ws.Range("E52").Value = 500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 139500
$ws.Range("E54").Value = 148900
$ws.Range("D57").Value = 5700
$ws.Range("E57").Value = 4900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 20400
$ws.Range("E59").Value = 19200
$ws.Range("D60").Value = 26200
$ws.Range("E60").Value = 24100
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 800
$ws.Range("E62").Value = 800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 27000
$ws.Range("E66").Value = 24900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -140500
$ws.Range("E72").Value = -126600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 112500
$ws.Range("E76").Value = 124000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -13900
$ws.Range("E81").Value = -8700
$ws.Range("D83").Value = 100
$ws.Range("E83").Value = 200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -10100
$ws.Range("E89").Value = -9300
$ws.Range("D91").Value = -100
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 4500
$ws.Range("E94").Value = 9900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 500
$ws.Range("E100").Value = 23000
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = -5200
$ws.Range("E102").Value = 23600
